# Updates cryptos price/volume figures in the active worksheet to match
# the latest scrape (GitHub Actions automated refresh).
# Cells in column D that look like plain decimal numbers are forced to
# stay as text (NumberFormat "@") before assignment so Excel does not
# silently reinterpret values such as "1.003" as the number 1.003.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.752.81"
$ws.Range("D3").Value = "1.737.96"
$ws.Range("E3").Value = "  +5.14%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.79"
$ws.Range("E5").Value = "  +4.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5456"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2760"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("E9").Value = "  +5.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.91"
$ws.Range("E10").Value = "  +6.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07777"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.695"
$ws.Range("E12").Value = "  +2.07%  "
$ws.Range("D13").Value = "1.731.02"
$ws.Range("E13").Value = "  +4.58%  "
$ws.Range("D14").Value = "1.976.40"
$ws.Range("E14").Value = "  +5.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5993"
$ws.Range("E15").Value = "  +6.74%  "
$ws.Range("D16").Value = "0.0₅8426"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.45"
$ws.Range("E17").Value = "  +5.91%  "
$ws.Range("D18").Value = "27.764.25"
$ws.Range("E18").Value = "  +6.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "226.48"
$ws.Range("E19").Value = "  +18.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.841"
$ws.Range("E20").Value = "  +3.43%  "
$ws.Range("E22").Value = "  +5.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.235"
$ws.Range("E23").Value = "  +4.44%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.18"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1254"
$ws.Range("E26").Value = "  +4.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.704"
$ws.Range("E27").Value = "  +12.19%  "
$ws.Range("E28").Value = "  +7.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.457"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05674"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.315"
$ws.Range("E31").Value = "  +3.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.699"
$ws.Range("E32").Value = "  +5.92%  "
$ws.Range("E33").Value = "  +4.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.685"
$ws.Range("E34").Value = "  +6.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9770"
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.854"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("E37").Value = "  +1.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5970"
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01671"
$ws.Range("E39").Value = "  +4.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.916"
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "1.051.79"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8491"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.13"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").Value = "1.881.08"
$ws.Range("E45").Value = "  +5.00%  "
$ws.Range("E46").Value = "  +14.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "59.61"
$ws.Range("E47").Value = "  +2.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.278"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05322"
$ws.Range("E51").Value = "  -0.15%  "
